# Apply updated classification results / explanations / scores
# per the commit: "Add updated notebooks and classification results, new interest rate data"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'The article suggests a hawkish tone as it emphasizes the need to combat inflation and return it to the 2% target. It highlights that inflation is being pushed down due to restrictive financing conditions and past interest rate increases. Additionally, the Governing Council mentions the possibility of reducing monetary policy restriction if confident in inflation convergence but also states it will continue to follow a data-dependent approach. Overall, the strong focus on price stability and willingness to adjust instruments to achieve the inflation target leans towards a hawkish stance.'
$ws.Range("F3").Value = '''0.6'
$ws.Range("G3").Value = 'The article indicates a hawkish stance as the ECB is focused on combating low inflation by maintaining restrictive policy rates for as long as necessary. The projections show a downward revision in inflation rates, with the Governing Council determined to ensure inflation returns to its 2% target. The emphasis on maintaining restrictive interest rates and reducing the PEPP portfolio suggests a hawkish sentiment towards monetary policy.'
$ws.Range("G4").Value = 'The article suggests a hawkish stance as it emphasizes the Governing Council''s determination to ensure that inflation returns to its 2% medium-term target in a timely manner by maintaining interest rates at restrictive levels. It highlights the impact of tight financing conditions on dampening demand and pushing down inflation, indicating a bias towards combating inflation rather than supporting economic growth.'
$ws.Range("G5").Value = 'The article suggests a hawkish sentiment by emphasizing the Governing Council''s determination to ensure inflation returns to the 2% target, maintaining interest rates at restrictive levels for as long as necessary, and advancing the normalisation of the balance sheet. The downward revision of inflation projections for 2023 and 2024, coupled with the focus on combating inflation and tightening monetary policy, contribute to the hawkish classification and score.'
$ws.Range("G6").Value = 'The article suggests a hawkish stance by indicating that inflation is expected to stay high, domestic price pressures are strong, and past interest rate increases are dampening demand to help push down inflation. The Governing Council is determined to ensure inflation returns to its 2% target, and future decisions will ensure that policy rates will be set at restrictive levels for as long as necessary. The emphasis on maintaining interest rates at levels that will contribute to reducing inflation and the readiness to adjust all instruments within the mandate to achieve the 2% target indicate a hawkish sentiment.'
$ws.Range("F7").Value = '''0.9'
$ws.Range("G7").Value = 'The article suggests a hawkish stance by highlighting the increase in ECB interest rates in response to high inflation and tightening financing conditions. The Governing Council is determined to combat inflation and bring it back to the 2% target. This indicates a prioritization of controlling inflation over supporting economic growth, which is a hawkish approach.'
$ws.Range("G8").Value = 'The article suggests a hawkish stance as the ECB has decided to raise the key interest rates by 25 basis points in order to combat inflation and bring it back to the 2% target. The language used in the article emphasizes the need for restrictive monetary policy for an extended period to achieve the inflation target, indicating a hawkish tone.'
$ws.Range("F9").Value = '''0.7'
$ws.Range("G9").Value = 'The article suggests a hawkish stance as the Governing Council decided to raise interest rates due to concerns about high inflation levels. The language used emphasizes the need to combat inflation and achieve the 2% medium-term target, indicating a hawkish sentiment towards monetary policy.'
$ws.Range("G10").Value = 'The article suggests hawkish sentiments as the European Central Bank (ECB) has decided to raise interest rates by 25 basis points to combat high inflation. The ECB aims to bring inflation back to its 2% target and indicates a willingness to implement restrictive policies as necessary to achieve this goal. This indicates a hawkish stance towards monetary policy.'
$ws.Range("E11").Value = 'Hawkish'
$ws.Range("F11").Value = '''0.8'
$ws.Range("G11").Value = 'The article indicates a hawkish stance by highlighting concerns about high inflation and announcing a decision to raise interest rates in order to combat inflation and bring it back to the 2% target. The language used regarding future policy decisions and the emphasis on tightening monetary policy support the hawkish classification.'
$ws.Range("G12").Value = 'The article suggests a hawkish stance by indicating that the ECB has decided to increase interest rates to combat high inflation. It emphasizes the importance of returning inflation to the 2% target and mentions the tightening of monetary policy to dampen demand pressures. Additionally, the article discusses resilience in the banking sector and the readiness to adjust instruments to ensure price stability, all of which lean towards a hawkish outlook.'
$ws.Range("G13").Value = 'The article suggests a hawkish stance by indicating the intention of raising interest rates significantly and keeping them at restrictive levels to combat inflation. The Governing Council aims to reduce inflation by dampening demand and guard against the risk of persistent increase in inflation expectations, which are typical hawkish policies. Additionally, the decision to reduce the Eurosystem''s holdings of securities under the asset purchase programme also reflects a hawkish tone towards monetary policy. The overall tone is focused on combating inflation and maintaining price stability.'
$ws.Range("F14").Value = '''0.8'
$ws.Range("G14").Value = 'The article suggests a hawkish stance as it highlights the decision to raise interest rates significantly to combat high inflation levels. The Governing Council expects to continue raising rates in a steady manner to reach restrictive levels. Additionally, the article mentions a decline in the asset purchase program and a focus on reducing inflation, which align with a hawkish monetary policy stance aimed at controlling rising prices.'
$ws.Range("G15").Value = 'The article clearly indicates a hawkish stance as the ECB has decided to raise interest rates by a significant amount, recognizing high inflation levels as a major concern. The emphasis on normalizing monetary policy to combat rising inflation pressures and the commitment to adjusting all instruments within its mandate to stabilize inflation at the 2% target suggest a hawkish tone. The score of 0.8 reflects the strong hawkish sentiment conveyed in the article.'
$ws.Range("G16").Value = 'The article suggests a hawkish stance as it describes the European Central Bank''s decision to raise interest rates multiple times to combat high inflation, with a focus on returning inflation to its 2% target. The emphasis on reducing support for demand, adjusting interest rates, and normalizing the monetary policy reflects a hawkish outlook towards controlling inflationary pressures.'
$ws.Range("F17").Value = '''0.8'
$ws.Range("G17").Value = 'The article strongly indicates a hawkish stance by the ECB as it mentions the decision to raise interest rates significantly to combat high inflation levels. The emphasis on further interest rate hikes and the expectation for inflation to remain above target for an extended period all point towards a hawkish outlook. Additionally, the lowered economic growth projections and concerns about geopolitical situations further support the hawkish classification.'
$ws.Range("E18").Value = 'Error'
$ws.Range("F18").Value = 'Error'
$ws.Range("G18").Value = 'Error'
$ws.Range("E19").Value = 'Error'
$ws.Range("F19").Value = 'Error'
$ws.Range("G19").Value = 'Error'
$ws.Range("G20").Value = 'The article suggests a hawkish stance by detailing the Governing Council''s decision to raise interest rates significantly, implement the Transmission Protection Instrument (TPI) to support monetary policy transmission, and signal further normalization of interest rates in the future. These actions are aimed at combating inflation and ensuring price stability in the medium term.'
$ws.Range("F21").Value = '''0.8'
$ws.Range("G21").Value = 'The article indicates a hawkish stance by discussing the challenges of high inflation and the measures being taken to combat it, such as raising interest rates and ending net asset purchases. The planned rate hikes and commitment to maintaining price stability signal a hawkish approach towards monetary policy.'
$ws.Range("E22").Value = 'Dovish'
$ws.Range("F22").Value = '''-0.7'
$ws.Range("G22").Value = 'The article suggests that the economy is facing challenges due to the conflict in Ukraine and the associated uncertainties. It highlights that economic activity is still being supported by the reopening after the pandemic crisis, but inflation has increased significantly mainly due to rising energy costs. The European Central Bank (ECB) is maintaining a flexible approach to monetary policy, with a focus on price stability. The overall tone leans towards dovishness as the ECB is prepared to adjust its policies as needed to support economic recovery and ensure price stability.'
$ws.Range("F23").Value = '''-0.5'
$ws.Range("G23").Value = 'The article suggests a dovish stance by emphasizing support for liquidity conditions, implementing sanctions for stability, and adjusting asset purchase programs to address the uncertain economic environment. The ECB maintains a commitment to maintaining accommodative monetary policy and supporting economic recovery. The dovish score of -0.5 reflects a moderately supportive tone towards economic growth and price stability.'
$ws.Range("G24").Value = 'The article suggests a dovish stance by detailing the continuation of asset purchase programs, the maintenance of low interest rates, and the commitment to maintaining accommodative monetary policy until inflation reaches the target of 2%. The flexibility and reinforcement of asset purchases indicate a supportive approach to economic growth and price stability. The score is closer to -1, indicating a strongly dovish sentiment.'
$ws.Range("F25").Value = '''0.6'
$ws.Range("G25").Value = 'The article suggests a hawkish stance by outlining a reduction in asset purchases, maintaining accommodative monetary policy, and signaling a possible increase in interest rates in the future. The Governing Council''s focus on inflation reaching and stabilizing at 2% and the intention to adjust instruments to ensure this target indicates a hawkish bias towards combating inflation.'
$ws.Range("F26").Value = '''-0.8'
$ws.Range("G26").Value = 'The article suggests a dovish stance as it emphasizes maintaining accommodative monetary policy measures to support economic recovery and achieve the inflation target of two per cent over the medium term. The continuation of net asset purchases under the PEPP, the unchanged interest rates, and the commitment to provide ample liquidity all reflect a supportive stance towards economic growth and low-interest rates.'
$ws.Range("G27").Value = 'The article suggests a continuation of accommodative monetary policy measures such as maintaining low interest rates, ongoing asset purchases, and ensuring ample liquidity in the financial system. The focus is on supporting the economy and maintaining favorable financing conditions, indicating a dovish stance overall. The score is slightly below neutral indicating a modestly dovish sentiment.'
$ws.Range("G28").Value = 'The article suggests a highly dovish stance as the European Central Bank (ECB) is committed to maintaining an accommodative monetary policy stance to meet its inflation target of two per cent. The ECB plans to keep interest rates low until inflation reaches the target well ahead of the projection horizon. The article also outlines continued asset purchases and liquidity provision measures to support economic recovery and price stability.'
$ws.Range("G29").Value = 'The article highlights a very accommodative monetary policy stance with unchanged low interest rates, continued net asset purchases, and a commitment to provide ample liquidity. The focus on supporting financing conditions and countering the downward impact of the pandemic on inflation suggests a dovish stance aimed at fostering economic growth and maintaining low interest rates.'
$ws.Range("G30").Value = 'The article suggests a highly accommodative monetary policy stance by keeping interest rates at historic lows and continuing with net asset purchases. The emphasis on maintaining favorable financing conditions and supporting the smooth transmission of monetary policy indicate a dovish stance aimed at countering the negative impact of the pandemic on inflation and economic growth. The commitment to adjusting all instruments as necessary to ensure inflation moves towards the target in a sustained manner also reflects a dovish predisposition.'
$ws.Range("G31").Value = 'The article suggests a dovish stance by indicating continued net asset purchases, maintaining low interest rates, and providing ample liquidity. The emphasis on flexibility in asset purchases and the commitment to maintaining accommodative policies until inflation reaches the desired level all point towards supportive measures for economic growth. The score of -0.8 reflects a moderately dovish stance.'
$ws.Range("G32").Value = 'The article suggests a very accommodative monetary policy stance with a focus on supporting economic growth and maintaining low interest rates. The continuation of asset purchases and reinvestment of principal payments indicate a commitment to providing ample liquidity and monetary accommodation. The language used in the article leans heavily towards dovish policies aimed at bolstering the economy and ensuring inflation converges to target levels.'
